$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 14 more match rows (16-29) below the existing 14 (rows 2-15) --
# same columns: venue, date, result, ownTeam, oppTeam, batsman, totalRuns,
# totalBalls, total4s, total6s, sr.
#
# Columns G:K hold numeric-looking values (run/ball/four/six counts, strike
# rate) that must stay stored as TEXT, matching every existing data cell on
# this sheet (plain "8" / "114.28", not the number 8 / 114.28). Force Text
# format on those columns before assigning so Excel does not auto-convert the
# strings to numbers, then drop back to the workbook default style afterwards
# so the new cells are not left with a stray explicit format (same unstyled
# look as the rest of the sheet).
$ws.Range("G16:K29").NumberFormat = "@"

# Row 16
$ws.Range("A16").Value = " Abu Dhabi"
$ws.Range("B16").Value = " November 02 2020"
$ws.Range("C16").Value = "Capitals won by 6 wickets (with 6 balls remaining)"
$ws.Range("D16").Value = "Delhi Capitals"
$ws.Range("E16").Value = "Royal Challengers Bangalore"
$ws.Range("F16").Value = "Rishabh Pant †"
$ws.Range("G16").Value = "8"
$ws.Range("H16").Value = "7"
$ws.Range("I16").Value = "1"
$ws.Range("J16").Value = "0"
$ws.Range("K16").Value = "114.28"

# Row 17
$ws.Range("A17").Value = " Dubai (DSC)"
$ws.Range("B17").Value = " November 05 2020"
$ws.Range("C17").Value = "Mumbai won by 57 runs"
$ws.Range("D17").Value = "Delhi Capitals"
$ws.Range("E17").Value = "Mumbai Indians"
$ws.Range("F17").Value = "Rishabh Pant †"
$ws.Range("G17").Value = "3"
$ws.Range("H17").Value = "9"
$ws.Range("I17").Value = "0"
$ws.Range("J17").Value = "0"
$ws.Range("K17").Value = "33.33"

# Row 18
$ws.Range("A18").Value = " Dubai (DSC)"
$ws.Range("B18").Value = " October 27 2020"
$ws.Range("C18").Value = "Sunrisers won by 88 runs"
$ws.Range("D18").Value = "Delhi Capitals"
$ws.Range("E18").Value = "Sunrisers Hyderabad"
$ws.Range("F18").Value = "Rishabh Pant †"
$ws.Range("G18").Value = "36"
$ws.Range("H18").Value = "35"
$ws.Range("I18").Value = "3"
$ws.Range("J18").Value = "1"
$ws.Range("K18").Value = "102.85"

# Row 19
$ws.Range("A19").Value = " Abu Dhabi"
$ws.Range("B19").Value = " September 29 2020"
$ws.Range("C19").Value = "Sunrisers won by 15 runs"
$ws.Range("D19").Value = "Delhi Capitals"
$ws.Range("E19").Value = "Sunrisers Hyderabad"
$ws.Range("F19").Value = "Rishabh Pant †"
$ws.Range("G19").Value = "28"
$ws.Range("H19").Value = "27"
$ws.Range("I19").Value = "1"
$ws.Range("J19").Value = "2"
$ws.Range("K19").Value = "103.70"

# Row 20
$ws.Range("A20").Value = " Abu Dhabi"
$ws.Range("B20").Value = " October 24 2020"
$ws.Range("C20").Value = "KKR won by 59 runs"
$ws.Range("D20").Value = "Delhi Capitals"
$ws.Range("E20").Value = "Kolkata Knight Riders"
$ws.Range("F20").Value = "Rishabh Pant †"
$ws.Range("G20").Value = "27"
$ws.Range("H20").Value = "33"
$ws.Range("I20").Value = "2"
$ws.Range("J20").Value = "1"
$ws.Range("K20").Value = "81.81"

# Row 21
$ws.Range("A21").Value = " Dubai (DSC)"
$ws.Range("B21").Value = " October 31 2020"
$ws.Range("C21").Value = "Mumbai won by 9 wickets (with 34 balls remaining)"
$ws.Range("D21").Value = "Delhi Capitals"
$ws.Range("E21").Value = "Mumbai Indians"
$ws.Range("F21").Value = "Rishabh Pant †"
$ws.Range("G21").Value = "21"
$ws.Range("H21").Value = "24"
$ws.Range("I21").Value = "2"
$ws.Range("J21").Value = "0"
$ws.Range("K21").Value = "87.50"

# Row 22
$ws.Range("A22").Value = " Dubai (DSC)"
$ws.Range("B22").Value = " November 10 2020"
$ws.Range("C22").Value = "Mumbai won by 5 wickets (with 8 balls remaining)"
$ws.Range("D22").Value = "Delhi Capitals"
$ws.Range("E22").Value = "Mumbai Indians"
$ws.Range("F22").Value = "Rishabh Pant †"
$ws.Range("G22").Value = "56"
$ws.Range("H22").Value = "38"
$ws.Range("I22").Value = "4"
$ws.Range("J22").Value = "2"
$ws.Range("K22").Value = "147.36"

# Row 23
$ws.Range("A23").Value = " Dubai (DSC)"
$ws.Range("B23").Value = " October 20 2020"
$ws.Range("C23").Value = "Kings XI won by 5 wickets (with 6 balls remaining)"
$ws.Range("D23").Value = "Delhi Capitals"
$ws.Range("E23").Value = "Kings XI Punjab"
$ws.Range("F23").Value = "Rishabh Pant †"
$ws.Range("G23").Value = "14"
$ws.Range("H23").Value = "20"
$ws.Range("I23").Value = "1"
$ws.Range("J23").Value = "0"
$ws.Range("K23").Value = "70.00"

# Row 24
$ws.Range("A24").Value = " Dubai (DSC)"
$ws.Range("B24").Value = " October 05 2020"
$ws.Range("C24").Value = "Capitals won by 59 runs"
$ws.Range("D24").Value = "Delhi Capitals"
$ws.Range("E24").Value = "Royal Challengers Bangalore"
$ws.Range("F24").Value = "Rishabh Pant †"
$ws.Range("G24").Value = "37"
$ws.Range("H24").Value = "25"
$ws.Range("I24").Value = "3"
$ws.Range("J24").Value = "2"
$ws.Range("K24").Value = "148.00"

# Row 25
$ws.Range("A25").Value = " Dubai (DSC)"
$ws.Range("B25").Value = " September 20 2020"
$ws.Range("C25").Value = "Match tied (Capitals won the one-over eliminator)"
$ws.Range("D25").Value = "Delhi Capitals"
$ws.Range("E25").Value = "Kings XI Punjab"
$ws.Range("F25").Value = "Rishabh Pant †"
$ws.Range("G25").Value = "31"
$ws.Range("H25").Value = "29"
$ws.Range("I25").Value = "4"
$ws.Range("J25").Value = "0"
$ws.Range("K25").Value = "106.89"

# Row 26
$ws.Range("A26").Value = " Abu Dhabi"
$ws.Range("B26").Value = " November 08 2020"
$ws.Range("C26").Value = "Capitals won by 17 runs"
$ws.Range("D26").Value = "Delhi Capitals"
$ws.Range("E26").Value = "Sunrisers Hyderabad"
$ws.Range("F26").Value = "Rishabh Pant †"
$ws.Range("G26").Value = "2"
$ws.Range("H26").Value = "3"
$ws.Range("I26").Value = "0"
$ws.Range("J26").Value = "0"
$ws.Range("K26").Value = "66.66"

# Row 27
$ws.Range("A27").Value = " Sharjah"
$ws.Range("B27").Value = " October 03 2020"
$ws.Range("C27").Value = "Capitals won by 18 runs"
$ws.Range("D27").Value = "Delhi Capitals"
$ws.Range("E27").Value = "Kolkata Knight Riders"
$ws.Range("F27").Value = "Rishabh Pant †"
$ws.Range("G27").Value = "38"
$ws.Range("H27").Value = "17"
$ws.Range("I27").Value = "5"
$ws.Range("J27").Value = "1"
$ws.Range("K27").Value = "223.52"

# Row 28
$ws.Range("A28").Value = " Dubai (DSC)"
$ws.Range("B28").Value = " September 25 2020"
$ws.Range("C28").Value = "Capitals won by 44 runs"
$ws.Range("D28").Value = "Delhi Capitals"
$ws.Range("E28").Value = "Chennai Super Kings"
$ws.Range("F28").Value = "Rishabh Pant †"
$ws.Range("G28").Value = "37"
$ws.Range("H28").Value = "25"
$ws.Range("I28").Value = "5"
$ws.Range("J28").Value = "0"
$ws.Range("K28").Value = "148.00"

# Row 29
$ws.Range("A29").Value = " Sharjah"
$ws.Range("B29").Value = " October 09 2020"
$ws.Range("C29").Value = "Capitals won by 46 runs"
$ws.Range("D29").Value = "Delhi Capitals"
$ws.Range("E29").Value = "Rajasthan Royals"
$ws.Range("F29").Value = "Rishabh Pant †"
$ws.Range("G29").Value = "5"
$ws.Range("H29").Value = "9"
$ws.Range("I29").Value = "0"
$ws.Range("J29").Value = "0"
$ws.Range("K29").Value = "55.55"

$ws.Range("G16:K29").Style = "Normal"
